$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the oneway-specific bike lane LTS rows (rows 64-68, where B=oneway='t')
$ws.Range("A64:G68").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)

# Scroll/select to match the saved view state
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 55
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C61").Select()
